$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), copying the existing header style from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data values for the Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
